# Apply the edits captured in the commit:
#  - "Checkout Your Information" sheet becomes the active/selected sheet
#    (tabSelected moves from "Login Details" to it, activeTab -> 1)
#  - selection on that sheet moves to L11
#  - cell C2 on that sheet changes from the number 1687 to the text "PN"
#    (adds a new shared string entry)

$wb = $excel.ActiveWorkbook

$wsCheckout = $wb.Worksheets.Item("Checkout Your Information")

# Change C2 from numeric 1687 to text "PN"
$wsCheckout.Range("C2").Value = "PN"

# Make "Checkout Your Information" the active sheet and select L11,
# which also clears tabSelected on the previously active "Login Details" sheet.
$wsCheckout.Activate()
$wsCheckout.Range("L11").Select()
